$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rowCount = $used.Rows.Count

for ($r = 2; $r -le $rowCount; $r++) {
    $cVal = $ws.Cells.Item($r, 3).Value()

    if ($cVal -eq "Running-Boys") {
        $ws.Cells.Item($r, 2).Value = "sports_club_boys"
        $ws.Cells.Item($r, 3).Value = "Running"
    }
    elseif ($cVal -eq "Running-Girls") {
        $ws.Cells.Item($r, 2).Value = "sports_club_girls"
        $ws.Cells.Item($r, 3).Value = "Running"
    }
    elseif ($cVal -eq "Running-Coed") {
        $ws.Cells.Item($r, 2).Value = "sports_club_coed"
        $ws.Cells.Item($r, 3).Value = "Running"
    }
}
